$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.5187906666666667
$ws.Cells.Item(2, 8).Value = 1.556372
$ws.Cells.Item(2, 9).Value = 0.5259328345599914
$ws.Cells.Item(2, 10).Value = 0.5259328345599914
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 63.255491791636
$ws.Cells.Item(2, 18).Value = 569.2994261247239
$ws.Cells.Item(2, 19).Value = 0.120030082456705
$ws.Cells.Item(2, 20).Value = 0.1272610425286561
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.5187906666666667
$ws.Cells.Item(3, 8).Value = 1.556372
$ws.Cells.Item(3, 9).Value = 0.5259328345599914
$ws.Cells.Item(3, 10).Value = 0.5259328345599914
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 76.73657387025334
$ws.Cells.Item(3, 18).Value = 690.62916483228
$ws.Cells.Item(3, 19).Value = 0.1456110296230348
$ws.Cells.Item(3, 20).Value = 0.154383060098142
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.5187906666666667
$ws.Cells.Item(4, 8).Value = 1.556372
$ws.Cells.Item(4, 9).Value = 0.5259328345599914
$ws.Cells.Item(4, 10).Value = 0.5259328345599914
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 43.32159871041956
$ws.Cells.Item(4, 18).Value = 389.894388393776
$ws.Cells.Item(4, 19).Value = 0.08220464212809275
$ws.Cells.Item(4, 20).Value = 0.08715688803837672
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.5187906666666667
$ws.Cells.Item(5, 8).Value = 1.556372
$ws.Cells.Item(5, 9).Value = 0.5259328345599914
$ws.Cells.Item(5, 10).Value = 0.5259328345599914
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 47.24546837268334
$ws.Cells.Item(5, 18).Value = 283.4728102361
$ws.Cells.Item(5, 19).Value = 0.08965035768212408
$ws.Cells.Item(5, 20).Value = 0.06336743671909226
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.5187906666666667
$ws.Cells.Item(6, 8).Value = 1.556372
$ws.Cells.Item(6, 9).Value = 0.5259328345599914
$ws.Cells.Item(6, 10).Value = 0.5259328345599914
$ws.Cells.Item(6, 13).Value = 89.83562999999999
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 46.60588637812
$ws.Cells.Item(6, 18).Value = 419.45297740308
$ws.Cells.Item(6, 19).Value = 0.08843672267003469
$ws.Cells.Item(6, 20).Value = 0.09376440717572429
$ws.Cells.Item(7, 7).Value = 0.4676293333333333
$ws.Cells.Item(7, 8).Value = 1.402888
$ws.Cells.Item(7, 9).Value = 0.4740671654400085
$ws.Cells.Item(7, 10).Value = 0.4740671654400086
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 57.01745493274399
$ws.Cells.Item(7, 18).Value = 513.1570943946959
$ws.Cells.Item(7, 19).Value = 0.10819313269419
$ws.Cells.Item(7, 20).Value = 0.1147110006032884
$ws.Cells.Item(8, 7).Value = 0.4676293333333333
$ws.Cells.Item(8, 8).Value = 1.402888
$ws.Cells.Item(8, 9).Value = 0.4740671654400085
$ws.Cells.Item(8, 10).Value = 0.4740671654400086
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 69.16907952834664
$ws.Cells.Item(8, 18).Value = 622.5217157551199
$ws.Cells.Item(8, 19).Value = 0.1312513757159599
$ws.Cells.Item(8, 20).Value = 0.1391583390185394
$ws.Cells.Item(9, 7).Value = 0.4676293333333333
$ws.Cells.Item(9, 8).Value = 1.402888
$ws.Cells.Item(9, 9).Value = 0.4740671654400085
$ws.Cells.Item(9, 10).Value = 0.4740671654400086
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 39.04937313936711
$ws.Cells.Item(9, 18).Value = 351.4443582543039
$ws.Cells.Item(9, 19).Value = 0.07409790588997729
$ws.Cells.Item(9, 20).Value = 0.07856177851206665
$ws.Cells.Item(10, 7).Value = 0.4676293333333333
$ws.Cells.Item(10, 8).Value = 1.402888
$ws.Cells.Item(10, 9).Value = 0.4740671654400085
$ws.Cells.Item(10, 10).Value = 0.4740671654400086
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 42.58628440656667
$ws.Cells.Item(10, 18).Value = 255.5177064394
$ws.Cells.Item(10, 19).Value = 0.08080935084154665
$ws.Cells.Item(10, 20).Value = 0.05711836024033707
$ws.Cells.Item(11, 7).Value = 0.4676293333333333
$ws.Cells.Item(11, 8).Value = 1.402888
$ws.Cells.Item(11, 9).Value = 0.4740671654400085
$ws.Cells.Item(11, 10).Value = 0.4740671654400086
$ws.Cells.Item(11, 13).Value = 89.83562999999999
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 42.00977576647999
$ws.Cells.Item(11, 18).Value = 378.08798189832
$ws.Cells.Item(11, 19).Value = 0.07971540029833458
$ws.Cells.Item(11, 20).Value = 0.084517687065777
